$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 17, shifting existing rows 17-24 down to 19-26.
$ws.Rows("17:18").Insert()

# New row 17 data (copy of A-K,Q,R from the (now shifted) row 19, with updated D,L,M,N,O,P,S)
$ws.Cells.Item(17,1).Value2 = 1
$ws.Cells.Item(17,2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(17,3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(17,4).Value2 = 45086
$ws.Cells.Item(17,5).Value2 = 15
$ws.Cells.Item(17,6).Value2 = "Fruta"
$ws.Cells.Item(17,7).Value2 = 100101
$ws.Cells.Item(17,8).Value2 = "Berries"
$ws.Cells.Item(17,9).Value2 = 100101007
$ws.Cells.Item(17,10).Value2 = "Kiwi"
$ws.Cells.Item(17,11).Value2 = "Hayward"
$ws.Cells.Item(17,12).Value2 = "Especial"
$ws.Cells.Item(17,13).Value2 = 250
$ws.Cells.Item(17,14).Value2 = 25000
$ws.Cells.Item(17,15).Value2 = 26000
$ws.Cells.Item(17,16).Value2 = 25500
$ws.Cells.Item(17,17).Value2 = "`$/bandeja 18 kilos"
$ws.Cells.Item(17,18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(17,19).Value2 = 1417
$ws.Cells.Item(17,20).Value2 = 18

# New row 18 data
$ws.Cells.Item(18,1).Value2 = 1
$ws.Cells.Item(18,2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(18,3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(18,4).Value2 = 45086
$ws.Cells.Item(18,5).Value2 = 15
$ws.Cells.Item(18,6).Value2 = "Fruta"
$ws.Cells.Item(18,7).Value2 = 100101
$ws.Cells.Item(18,8).Value2 = "Berries"
$ws.Cells.Item(18,9).Value2 = 100101007
$ws.Cells.Item(18,10).Value2 = "Kiwi"
$ws.Cells.Item(18,11).Value2 = "Hayward"
$ws.Cells.Item(18,12).Value2 = "Primera"
$ws.Cells.Item(18,13).Value2 = 250
$ws.Cells.Item(18,14).Value2 = 20000
$ws.Cells.Item(18,15).Value2 = 21000
$ws.Cells.Item(18,16).Value2 = 20500
$ws.Cells.Item(18,17).Value2 = "`$/bandeja 18 kilos"
$ws.Cells.Item(18,18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(18,19).Value2 = 1139
$ws.Cells.Item(18,20).Value2 = 18
